# dang nhap va dang xuat
$wb = $excel.ActiveWorkbook

# --- users sheet: rename "username" column to "name", rename role values,
#     and fix avatar asset paths (asset -> assets, with a typo on avatar-5) ---
$users = $wb.Worksheets.Item("users")

$users.Range("B1").Value = "name"

$users.Range("A2").Value = "assets/jpg/avatar-2.jpg"
$users.Range("D2").Value = "administrator"

$users.Range("A3").Value = "assets/jpg/avatar-1.jpg"
$users.Range("D3").Value = "contact"

$users.Range("A4").Value = "assets/jpg/avatar-4.jpg"
$users.Range("D4").Value = "teacher"

$users.Range("A5").Value = "assetss/jpg/avatar-5.jpg"
$users.Range("D5").Value = "administrator"

$users.Range("A6").Value = "assets/jpg/avatar-3.jpg"
$users.Range("D6").Value = "student"

$users.Range("A7").Value = "assets/jpg/avatar-3.jpg"
$users.Range("D7").Value = "student"

$users.Range("A8").Value = "assets/jpg/avatar-3.jpg"
$users.Range("D8").Value = "student"

# --- move the active tab from "loai_khoa_hoc" to "users", with the
#     selection left at F28 on the users sheet ---
$loaiKhoaHoc = $wb.Worksheets.Item("loai_khoa_hoc")
$loaiKhoaHoc.Select()
$loaiKhoaHoc.Range("C18").Select()

$users.Select()
$users.Range("F28").Select()

$users.Activate()
